$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887

# Sheet ALC, Row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 439552.4
$ws.Range("I64").Value = 836850.2
$ws.Range("K64").Value = 836850.2
$ws.Range("M64").Value = -836602.2

# Sheet ALC, Row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 439552.4
$ws.Range("I67").Value = 836850.2
$ws.Range("K67").Value = 836850.2
$ws.Range("M67").Value = -835992.2

# Sheet ALC, Row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1490.9395
$ws.Range("J129").Value = 1760.6154
$ws.Range("L129").Value = 5281.8462
$ws.Range("N129").Value = -15281.8462

# Sheet ALC, Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27028334
$ws.Range("I137").Value = 45455532
$ws.Range("J137").Value = 1778.5333
$ws.Range("K137").Value = 136366596
$ws.Range("L137").Value = 5335.5999
$ws.Range("M137").Value = -136364046
$ws.Range("N137").Value = -10435.5999

# Sheet ALC, Row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1884.9445
$ws.Range("I141").Value = 1201.5
$ws.Range("K141").Value = 3604.5
$ws.Range("M141").Value = 1575.5

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15056.39
$ws.Range("I32").Value = 1992.1781
$ws.Range("J32").Value = 253478.25
$ws.Range("K32").Value = 1992.1781
$ws.Range("L32").Value = 253478.25
$ws.Range("M32").Value = -1705.1781
$ws.Range("N32").Value = -254052.25

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2466.9348
$ws.Range("I61").Value = 1571.1613
$ws.Range("J61").Value = 4318.2
$ws.Range("K61").Value = 1571.1613
$ws.Range("L61").Value = 4318.2
$ws.Range("M61").Value = -1359.1613
$ws.Range("N61").Value = -4742.2

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1678.3062
$ws.Range("I132").Value = 1355.2222
$ws.Range("K132").Value = 4065.6666
$ws.Range("M132").Value = -1535.6666

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2466.9348
$ws.Range("I136").Value = 1571.1613
$ws.Range("J136").Value = 4318.2
$ws.Range("K136").Value = 4713.4839
$ws.Range("L136").Value = 12954.6
$ws.Range("M136").Value = -2163.4839
$ws.Range("N136").Value = -18054.6

# Sheet ARM, Row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 47930.8
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47930.8
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47930.8
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -58210.8

# Sheet BSM, Row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2904.2632
$ws.Range("I20").Value = 2784.3572
$ws.Range("J20").Value = 3240
$ws.Range("K20").Value = 2784.3572
$ws.Range("L20").Value = 3240
$ws.Range("M20").Value = -2537.3572
$ws.Range("N20").Value = -3734

# Sheet CRP, Row 5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2862.5557
$ws.Range("I5").Value = 249.25
$ws.Range("K5").Value = 249.25
$ws.Range("M5").Value = -137.25

# Sheet CRP, Row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# Sheet CRP, Row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 726.5
$ws.Range("I19").Value = 805
$ws.Range("J19").Value = 700.3333
$ws.Range("K19").Value = 805
$ws.Range("L19").Value = 700.3333
$ws.Range("M19").Value = -635
$ws.Range("N19").Value = -1040.3333

# Sheet CRP, Row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 552.7273
$ws.Range("I22").Value = 280.16666
$ws.Range("K22").Value = 280.16666
$ws.Range("M22").Value = 69.83334000000002

# Sheet CRP, Row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 726.5
$ws.Range("I24").Value = 805
$ws.Range("J24").Value = 700.3333
$ws.Range("K24").Value = 805
$ws.Range("L24").Value = 700.3333
$ws.Range("M24").Value = -635
$ws.Range("N24").Value = -1040.3333

# Sheet CRP, Row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 62500000
$ws.Range("I99").Value = 62500000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 62500000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -62498502
$ws.Range("N99").ClearContents()

# Sheet CRP, Row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 62500000
$ws.Range("I126").Value = 62500000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 187500000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -187497530
$ws.Range("N126").ClearContents()

# Sheet CRP, Row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 28153.037
$ws.Range("J133").Value = 28455.23
$ws.Range("L133").Value = 28455.23
$ws.Range("N133").Value = -33515.23

# Sheet CRP, Row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 39420.715
$ws.Range("J135").Value = 39420.715
$ws.Range("L135").Value = 39420.715
$ws.Range("N135").Value = -49560.715

# Sheet CUL, Row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 582.48
$ws.Range("I122").Value = 257.70587
$ws.Range("J122").Value = 1272.625
$ws.Range("K122").Value = 2319.35283
$ws.Range("L122").Value = 11453.625
$ws.Range("M122").Value = 130.6471700000002
$ws.Range("N122").Value = -16353.625

# Sheet GSM, Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 10006
$ws.Range("I102").Value = 10012
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 10012
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -8390
$ws.Range("N102").Value = -13244

# Sheet GSM, Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2708.913
$ws.Range("I132").Value = 2490.725
$ws.Range("K132").Value = 7472.174999999999
$ws.Range("M132").Value = -4942.174999999999

# Sheet LTW, Row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1062.579
$ws.Range("I93").Value = 403.16666
$ws.Range("J93").Value = 2193
$ws.Range("K93").Value = 403.16666
$ws.Range("L93").Value = 2193
$ws.Range("M93").Value = 844.83334
$ws.Range("N93").Value = -4689

# Sheet LTW, Row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3375.585
$ws.Range("I136").Value = 2149.3809
$ws.Range("J136").Value = 8057.4546
$ws.Range("K136").Value = 6448.1427
$ws.Range("L136").Value = 24172.3638
$ws.Range("M136").Value = -3898.1427
$ws.Range("N136").Value = -29272.3638

# Sheet WVR, Row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 113201.445
$ws.Range("I122").Value = 251435.75
$ws.Range("J122").Value = 2614
$ws.Range("K122").Value = 754307.25
$ws.Range("L122").Value = 7842
$ws.Range("M122").Value = -751857.25
$ws.Range("N122").Value = -12742

# Sheet WVR, Row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 101310.6
$ws.Range("I126").Value = 250650.25
$ws.Range("K126").Value = 751950.75
$ws.Range("M126").Value = -749480.75

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8573637
$ws.Range("I136").Value = 9036975
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 27110925
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -27108375
$ws.Range("N136").Value = -10800
